# Backup before dimension reduction:
# Shift the "q" labels in column A (rows 2..97) down by one index,
# i.e. A2 "q1" -> "q0", A3 "q2" -> "q1", ..., A97 "q96" -> "q95".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 97; $r++) {
    $newIndex = $r - 2
    $ws.Cells.Item($r, 1).Value = "q$newIndex"
}
